# tdf#108064 OOXML export: keep preset dashes with any line width
#
# This edit touches three distinct, low-level pieces of OOXML that are not
# reachable through the high-level Word object model (a DrawingML custDash
# list inside a VML/DrawingML fallback shape, the w14:anchorId of the VML
# fallback <v:line>, and a handful of w:rsid attributes inside the Header /
# Footer style definitions). We therefore round-trip the whole package
# through Content.WordOpenXML, perform precise, scoped string edits, and
# write the result back.

$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

# ---------------------------------------------------------------------
# 1) Duplicate the last <a:ds/> entry in the custDash list so there are
#    four dash/space pairs instead of three.
# ---------------------------------------------------------------------
$custDashOld = '<a:ds d="800000" sp="300000"/><a:ds d="100000" sp="300000"/><a:ds d="100000" sp="300000"/></a:custDash>'
$custDashNew = '<a:ds d="800000" sp="300000"/><a:ds d="100000" sp="300000"/><a:ds d="100000" sp="300000"/><a:ds d="100000" sp="300000"/></a:custDash>'
if ($xml.IndexOf($custDashOld) -lt 0) {
    throw "custDash target not found"
}
$xml = $xml.Replace($custDashOld, $custDashNew)

# ---------------------------------------------------------------------
# 2) Re-stamp the w14:anchorId of the VML fallback <v:line> connector.
# ---------------------------------------------------------------------
$anchorOld = 'w14:anchorId="707FA7C3"'
$anchorNew = 'w14:anchorId="768003F5"'
if ($xml.IndexOf($anchorOld) -lt 0) {
    throw "anchorId target not found"
}
$xml = $xml.Replace($anchorOld, $anchorNew)

# ---------------------------------------------------------------------
# 3) Update the w:rsid stamps on the Header / HeaderChar / Footer /
#    FooterChar style definitions (word/styles.xml) only -- the same
#    rsid value also appears (and must stay untouched) in other parts
#    such as settings.xml, headers, footers, footnotes and endnotes, so
#    scope the replacement to the styles.xml package part.
# ---------------------------------------------------------------------
$partStartTag = '<pkg:part pkg:name="/word/styles.xml"'
$partStartIdx = $xml.IndexOf($partStartTag)
if ($partStartIdx -lt 0) {
    throw "styles.xml part not found"
}
$partEndTag = "</pkg:part>"
$partEndIdx = $xml.IndexOf($partEndTag, $partStartIdx) + $partEndTag.Length

$stylesPart = $xml.Substring($partStartIdx, $partEndIdx - $partStartIdx)
$rsidOld = 'w:rsid w:val="008C04E2"'
$rsidNew = 'w:rsid w:val="00DD09CB"'
$rsidCount = ([regex]::Matches($stylesPart, [regex]::Escape($rsidOld))).Count
if ($rsidCount -ne 4) {
    throw ("expected 4 rsid occurrences in styles.xml, found " + $rsidCount)
}
$stylesPartNew = $stylesPart.Replace($rsidOld, $rsidNew)

$xml = $xml.Substring(0, $partStartIdx) + $stylesPartNew + $xml.Substring($partEndIdx)

# ---------------------------------------------------------------------
# Write the modified package XML back to the document.
# ---------------------------------------------------------------------
$d.Content.WordOpenXML = $xml

Write-Output "applied"
